# Auto-generated: update Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(2, 4) "31.033.19"
Set-TextValue $ws.Cells.Item(2, 5) "  +0.86%  "

Set-TextValue $ws.Cells.Item(3, 4) "1.959.81"
Set-TextValue $ws.Cells.Item(3, 5) "  -0.12%  "

Set-TextValue $ws.Cells.Item(4, 5) "  +0.13%  "

Set-TextValue $ws.Cells.Item(5, 4) "244.57"
Set-TextValue $ws.Cells.Item(5, 5) "  -1.68%  "

Set-TextValue $ws.Cells.Item(6, 5) "  -0.16%  "

Set-TextValue $ws.Cells.Item(7, 4) "0.4875"
Set-TextValue $ws.Cells.Item(7, 5) "  +0.83%  "

Set-TextValue $ws.Cells.Item(8, 4) "0.2952"
Set-TextValue $ws.Cells.Item(8, 5) "  +0.67%  "

Set-TextValue $ws.Cells.Item(9, 4) "0.06987"
Set-TextValue $ws.Cells.Item(9, 5) "  +3.36%  "

Set-TextValue $ws.Cells.Item(10, 4) "19.47"
Set-TextValue $ws.Cells.Item(10, 5) "  +2.18%  "

Set-TextValue $ws.Cells.Item(11, 4) "107.93"
Set-TextValue $ws.Cells.Item(11, 5) "  -1.09%  "

Set-TextValue $ws.Cells.Item(12, 4) "1.951.11"
Set-TextValue $ws.Cells.Item(12, 5) "  -0.63%  "

Set-TextValue $ws.Cells.Item(13, 4) "0.07803"
Set-TextValue $ws.Cells.Item(13, 5) "  +0.80%  "

Set-TextValue $ws.Cells.Item(14, 4) "5.501"
Set-TextValue $ws.Cells.Item(14, 5) "  +0.70%  "

Set-TextValue $ws.Cells.Item(15, 4) "0.7014"
Set-TextValue $ws.Cells.Item(15, 5) "  +0.54%  "

Set-TextValue $ws.Cells.Item(16, 4) "281.19"
Set-TextValue $ws.Cells.Item(16, 5) "  -3.81%  "

Set-TextValue $ws.Cells.Item(17, 4) "31.045.95"
Set-TextValue $ws.Cells.Item(17, 5) "  +0.76%  "

Set-TextValue $ws.Cells.Item(18, 5) "  +1.35%  "

Set-TextValue $ws.Cells.Item(19, 4) "0.000007785"
Set-TextValue $ws.Cells.Item(19, 5) "  +1.12%  "

Set-TextValue $ws.Cells.Item(20, 4) "2.243.09"
Set-TextValue $ws.Cells.Item(20, 5) "  +1.10%  "

Set-TextValue $ws.Cells.Item(21, 5) "  -0.26%  "

Set-TextValue $ws.Cells.Item(22, 4) "5.544"
Set-TextValue $ws.Cells.Item(22, 5) "  -2.12%  "

Set-TextValue $ws.Cells.Item(23, 5) "  +0.17%  "

Set-TextValue $ws.Cells.Item(24, 4) "6.535"
Set-TextValue $ws.Cells.Item(24, 5) "  -1.11%  "

Set-TextValue $ws.Cells.Item(25, 4) "9.866"
Set-TextValue $ws.Cells.Item(25, 5) "  -0.29%  "

Set-TextValue $ws.Cells.Item(26, 4) "168.54"
Set-TextValue $ws.Cells.Item(26, 5) "  -1.04%  "

Set-TextValue $ws.Cells.Item(27, 4) "19.97"
Set-TextValue $ws.Cells.Item(27, 5) "  -0.42%  "

Set-TextValue $ws.Cells.Item(28, 4) "2.194"
Set-TextValue $ws.Cells.Item(28, 5) "  +0.88%  "

Set-TextValue $ws.Cells.Item(29, 5) "  -2.22%  "

Set-TextValue $ws.Cells.Item(30, 5) "  -4.09%  "

Set-TextValue $ws.Cells.Item(31, 4) "4.636"
Set-TextValue $ws.Cells.Item(31, 5) "  -4.75%  "

Set-TextValue $ws.Cells.Item(32, 5) "  -2.15%  "

Set-TextValue $ws.Cells.Item(33, 5) "  -0.12%  "

Set-TextValue $ws.Cells.Item(34, 4) "0.04925"
Set-TextValue $ws.Cells.Item(34, 5) "  -3.55%  "

Set-TextValue $ws.Cells.Item(35, 4) "0.7563"
Set-TextValue $ws.Cells.Item(35, 5) "  -1.93%  "

Set-TextValue $ws.Cells.Item(36, 5) "  -0.58%  "

Set-TextValue $ws.Cells.Item(37, 5) "  +0.18%  "

Set-TextValue $ws.Cells.Item(38, 5) "  -1.70%  "

Set-TextValue $ws.Cells.Item(39, 4) "2.707"
Set-TextValue $ws.Cells.Item(39, 5) "  -0.51%  "

Set-TextValue $ws.Cells.Item(40, 4) "6.568"
Set-TextValue $ws.Cells.Item(40, 5) "  +0.83%  "

Set-TextValue $ws.Cells.Item(41, 4) "78.14"
Set-TextValue $ws.Cells.Item(41, 5) "  +11.67%  "

Set-TextValue $ws.Cells.Item(42, 4) "2.137"
Set-TextValue $ws.Cells.Item(42, 5) "  +0.32%  "

Set-TextValue $ws.Cells.Item(43, 4) "0.9027"
Set-TextValue $ws.Cells.Item(43, 5) "  +1.41%  "

Set-TextValue $ws.Cells.Item(44, 4) "109.59"
Set-TextValue $ws.Cells.Item(44, 5) "  -0.48%  "

Set-TextValue $ws.Cells.Item(45, 4) "0.4456"
Set-TextValue $ws.Cells.Item(45, 5) "  -0.14%  "

Set-TextValue $ws.Cells.Item(46, 4) "8.139"
Set-TextValue $ws.Cells.Item(46, 5) "  +9.05%  "

Set-TextValue $ws.Cells.Item(47, 5) "  -0.15%  "

Set-TextValue $ws.Cells.Item(48, 4) "1.011.41"
Set-TextValue $ws.Cells.Item(48, 5) "  +8.46%  "

Set-TextValue $ws.Cells.Item(49, 4) "9.402"
Set-TextValue $ws.Cells.Item(49, 5) "  +0.38%  "

Set-TextValue $ws.Cells.Item(50, 5) "  -1.50%  "

Set-TextValue $ws.Cells.Item(51, 4) "36.03"
Set-TextValue $ws.Cells.Item(51, 5) "  +0.03%  "
